$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "roboticS1Prep" column (I) currently stores the text "No" for every
# data row (2-41). Convert it to a real boolean FALSE value, displayed via
# a custom "TRUE";"TRUE";"FALSE" number format, as part of cleaning up the
# database.

# Rows 2-27 already share one cell style (font Arial 10 black / General
# format), so a direct in-place edit converges on a single new style.
$g1 = $ws.Range("I2:I27")
$g1.Value = $false
$g1.NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Rows 28-41 started out on a different (default) cell style. Set their
# values first, then copy the already-normalized format from I2 onto them
# so every row ends up sharing the exact same resulting style instead of
# fragmenting into extra near-duplicate styles.
$g2 = $ws.Range("I28:I41")
$g2.Value = $false

$src = $ws.Range("I2")
$src.Copy()
$g2.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Reflect the edited column in the sheet's active selection/view, matching
# where the user's cursor ended up after making this pass through column I.
[void]$ws.Range("I2:I41").Select()
